$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.000.72'
$ws.Range("E2").Value = '  +3.56%  '
$ws.Range("D3").Value = '2.645.53'
$ws.Range("E3").Value = '  +5.86%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '113.78'
$ws.Range("E5").Value = '  +7.13%  '
$ws.Range("D6").Value = '326.56'
$ws.Range("E6").Value = '  +2.71%  '
$ws.Range("E7").Value = '  +2.07%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.556'
$ws.Range("E9").Value = '  +3.42%  '
$ws.Range("D10").Value = '40.86'
$ws.Range("E10").Value = '  +5.09%  '
$ws.Range("D11").Value = '20.14'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '0.0822'
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = '7.38'
$ws.Range("E14").Value = '  +3.90%  '
$ws.Range("D15").Value = '3.064.04'
$ws.Range("D16").Value = '2.643.59'
$ws.Range("E16").Value = '  +6.18%  '
$ws.Range("D17").Value = '0.872'
$ws.Range("E17").Value = '  +4.96%  '
$ws.Range("D18").Value = '49.967.48'
$ws.Range("E18").Value = '  +3.81%  '
$ws.Range("D19").Value = '13.22'
$ws.Range("E19").Value = '  +2.65%  '
$ws.Range("D20").Value = '6.76'
$ws.Range("E20").Value = '  +2.58%  '
$ws.Range("E21").Value = '  -3.06%  '
$ws.Range("D22").Value = '0.0₃0956'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = '72.24'
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("D24").Value = '275.84'
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("D26").Value = '26.83'
$ws.Range("E26").Value = '  +3.99%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  +3.00%  '
$ws.Range("E29").Value = '  +1.58%  '
$ws.Range("D30").Value = '36.40'
$ws.Range("E30").Value = '  +4.96%  '
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("D32").Value = '50.15'
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("E33").Value = '  +3.17%  '
$ws.Range("D34").Value = '19.55'
$ws.Range("E34").Value = '  +2.22%  '
$ws.Range("E35").Value = '  +5.39%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = '5.02'
$ws.Range("E37").Value = '  +9.37%  '
$ws.Range("E38").Value = '  +6.45%  '
$ws.Range("E39").Value = '  +7.83%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '124.06'
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.113'
$ws.Range("E41").Value = '  +1.88%  '
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").Value = '21.92'
$ws.Range("E43").Value = '  -3.23%  '
$ws.Range("E44").Value = '  +4.66%  '
$ws.Range("D45").Value = '2.087.10'
$ws.Range("E45").Value = '  +4.19%  '
$ws.Range("D46").Value = '3.33'
$ws.Range("E46").Value = '  +5.98%  '
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +13.53%  '
$ws.Range("D48").Value = '1.99'
$ws.Range("E48").Value = '  +4.53%  '
$ws.Range("D49").Value = '9.15'
$ws.Range("E49").Value = '  +2.42%  '
$ws.Range("E50").Value = '  +4.03%  '
$ws.Range("D51").Value = '59.60'
$ws.Range("E51").Value = '  +5.03%  '
